# Release23.06.01 Backlog Grooming Done
#
# 1) The cached "datetimeFigureOut" date field on the Slide Master and on
#    every one of the 11 slide layouts is bumped from 26/05/2023 to
#    29/05/2023 (the deck was re-saved three days later).
# 2) The "Hurry! Admissions filling up fast!!!" textbox on slide 1 is
#    dragged to a new position.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached Date Placeholder text everywhere it appears.
# ---------------------------------------------------------------------
$oldDate = "26/05/2023"
$newDate = "29/05/2023"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder($p.SlideMaster.Shapes)

# Every slide layout (CustomLayout) hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li).Shapes)
}

# ---------------------------------------------------------------------
# 2) Reposition the rotated "TextBox 33" promo textbox on slide 1.
# ---------------------------------------------------------------------
$EMU_PER_POINT = 12700
# Nudge by a hair so the point->EMU round trip lands on the exact target
# EMU value instead of truncating one EMU short.
$epsilon = 0.00001

$s1 = $p.Slides.Item(1)
$promo = $s1.Shapes.Item("TextBox 33")
$promo.Left = (7049767 / $EMU_PER_POINT) + $epsilon
$promo.Top = (1721632 / $EMU_PER_POINT) + $epsilon
